# Week 16 box score added a new pass-catcher (D.Sills) to the Giants'
# player list. Both the "Rushing" and "Receiving" trackers share the same
# column layout (one column per player), so insert a new column for him on
# each sheet -- right before the existing "E.Engram" column -- shifting the
# later columns over by one, then seed his header/placeholder data cell
# the same way every other player's column was initialised.
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new column before column R (currently "E.Engram"), shifting
    # the existing R:V columns right to S:W.
    $ws.Range("R1:R2").EntireColumn.Insert()

    # New player column: name header in row 1, placeholder "n" data in row 2
    # -- matching how every other player column is seeded before real yard
    # totals get filled in.
    $ws.Range("R1").Value = "D.Sills"
    $ws.Range("R2").Value = "n"
}
